# Applies the "habilidad" (skill) lookup-table normalization edit:
#  - L4 header renamed from "habilidad" to "id_habilidad"
#  - L5 value changed from the text "volar" to the numeric foreign key 1
#  - New small "habilidad" table added at F11:G13
#  - Selection moved to G14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the ninja table's "habilidad" column to "id_habilidad" and
# replace the literal skill name with a numeric reference id.
$ws.Range("L4").Value = "id_habilidad"
$ws.Range("L5").Value = 1

# New lookup table for habilidades (skills).
$ws.Range("F11").Value = "habilidad"
$ws.Range("F12").Value = "id"
$ws.Range("G12").Value = "habilidad"
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = "volar"

# Update the selection to match the recorded cursor position after the edit.
$ws.Range("G14").Select()
